$p = $ppt.ActivePresentation

# --- 1. Swap the table style on the table in slide 16 (Google Shape;213;p29)
#        from the custom "Table_0" style to the built-in themed table style.
$s16 = $p.Slides.Item(16)
for ($i = 1; $i -le $s16.Shapes.Count; $i++) {
    $sh = $s16.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{78EB1A86-E797-4C67-AAC9-3B746FCE418C}")
    }
}

# --- 2. Re-apply the "Office Theme" design (its colour scheme) to the deck,
#        replacing the current "Integral" theme colours.
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
